$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to align with new Excel file structure
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "IDAM Roles"

# Update the active selection to C1
$ws.Range("C1").Select()
